$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("展览").Range("F3").Value = 578
$wb.Worksheets.Item("展览").Range("F6").Value = 3177
$wb.Worksheets.Item("展览").Range("F7").Value = 2749
$wb.Worksheets.Item("展览").Range("F9").Value = 47
$wb.Worksheets.Item("展览").Range("F11").Value = 349
$wb.Worksheets.Item("展览").Range("F12").Value = 288
$wb.Worksheets.Item("展览").Range("F14").Value = 5707
$wb.Worksheets.Item("展览").Range("F16").Value = 1022
$wb.Worksheets.Item("展览").Range("F17").Value = 58
$wb.Worksheets.Item("展览").Range("F20").Value = 461
$wb.Worksheets.Item("展览").Range("F21").Value = 1241
$wb.Worksheets.Item("演出").Range("F9").Value = 60
$wb.Worksheets.Item("演出").Range("F24").Value = 291
$wb.Worksheets.Item("演出").Range("F25").Value = 4032
$wb.Worksheets.Item("演出").Range("F32").Value = 186
$wb.Worksheets.Item("本地生活").Range("F3").Value = 84
$wb.Worksheets.Item("本地生活").Range("F5").Value = 2565
$wb.Worksheets.Item("本地生活").Range("I5").Value = "//i2.hdslb.com/bfs/openplatform/202402/JAjBoc4t1708314351453.png"
$wb.Worksheets.Item("本地生活").Range("F9").Value = 1443
$wb.Worksheets.Item("全部类型").Range("F3").Value = 84
$wb.Worksheets.Item("全部类型").Range("F5").Value = 2565
$wb.Worksheets.Item("全部类型").Range("I5").Value = "//i2.hdslb.com/bfs/openplatform/202402/JAjBoc4t1708314351453.png"
$wb.Worksheets.Item("全部类型").Range("F7").Value = 1443
$wb.Worksheets.Item("全部类型").Range("F12").Value = 578
$wb.Worksheets.Item("全部类型").Range("F14").Value = 3177
$wb.Worksheets.Item("全部类型").Range("F15").Value = 2749
$wb.Worksheets.Item("全部类型").Range("F17").Value = 47
$wb.Worksheets.Item("全部类型").Range("F19").Value = 349
$wb.Worksheets.Item("全部类型").Range("F21").Value = 60
$wb.Worksheets.Item("全部类型").Range("F22").Value = 288
$wb.Worksheets.Item("全部类型").Range("F26").Value = 1022
$wb.Worksheets.Item("全部类型").Range("F28").Value = 58
$wb.Worksheets.Item("全部类型").Range("F31").Value = 461
$wb.Worksheets.Item("全部类型").Range("F38").Value = 291
$wb.Worksheets.Item("全部类型").Range("F39").Value = 1241
$wb.Worksheets.Item("全部类型").Range("F44").Value = 186
